# Update risk-driver calibration estimates (std moving average update).
$wb = $excel.ActiveWorkbook

# --- Sheet "Linear" ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.8033749316793535
$wsLinear.Range("B3").Value = 23.98165683982647
$wsLinear.Range("B4").Value = 31091.73189788084
$wsLinear.Range("B5").Value = "[1.0, 0.16405946202754063, 0.01992185860063315, 0.0008635681168387822, 0.004764298106171536, -0.050784052872094844, 0.10197142096963326, 0.22066897480243294, 0.07162675330864138, -0.04330517020340959, -0.03361440332166738, -0.03477593103295779, -0.06791268024199044, 0.09041901999835533, 0.19967688661231706, 0.042131623848433324, -0.050844573429099316, -0.02872307657753268, -0.017024791230251704, -0.04174144178641617]"

# --- Sheet "NonLinear" ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = -12.1669023512251
$wsNonLinear.Range("B5").Value = -1.527546859343692
$wsNonLinear.Range("B6").Value = 28780.35391094684
$wsNonLinear.Range("B7").Value = -10.79778203648345
$wsNonLinear.Range("B8").Value = 46.84307201573333
$wsNonLinear.Range("B9").Value = 33313.94580697861
$wsNonLinear.Range("B10").Value = "[0.9999999999999998, 0.1623576935663242, 0.02528190233573785, 0.0020282920221196253, 0.0029338757908369205, -0.04860922856308539, 0.09946308139233911, 0.21666984135919679, 0.07023577912713212, -0.04336960841539495, -0.03582757883235821, -0.03725052119477187, -0.06619565848038352, 0.08881236194204344, 0.19697301548465157, 0.04036484105107803, -0.05152728265936613, -0.02908167763294477, -0.015629050060204364, -0.0409388379107304]"
